$d = $word.ActiveDocument

# 1) "Auxiliar de coformacion" -> "desarrollador"
$d.Content.Find.Execute("Auxiliar de coformacion", $true, $false, $false, $false, $false,
                         $true, 1, $false, "desarrollador", 2)

# 2) "Mediante un contrato a Término Indefinido." -> "Mediante un contrato a Término Fijo."
$d.Content.Find.Execute("Mediante un contrato a Término Indefinido.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mediante un contrato a Término Fijo.", 2)

# 3) "Actualmente vigente desde el 2024-11-28." -> "Desde el 2024-01-22 hasta el 2024-01-30."
$d.Content.Find.Execute("Actualmente vigente desde el 2024-11-28.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Desde el 2024-01-22 hasta el 2024-01-30.", 2)

# 4) " devengando un salario de `$ 15000000." -> " devengando un salario de `$ 1000000."
$d.Content.Find.Execute(" devengando un salario de `$ 15000000.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " devengando un salario de `$ 1000000.", 2)

# 5) "  (01) días del mes de (febrero) de 2024" -> "  (05) días del mes de (febrero) de 2024"
$d.Content.Find.Execute("  (01) días del mes de (febrero) de 2024", $true, $false, $false, $false, $false,
                         $true, 1, $false, "  (05) días del mes de (febrero) de 2024", 2)
